$wb = $excel.ActiveWorkbook

# --- CMS sheet (sheet1) ---
$cms = $wb.Worksheets.Item("CMS")
$cms.Range("D25").Value = "WXX_XX_XXX_XXXBN0300"
$cms.Range("B26").Value = "WXX_XX_XXX_XXXBN0400"

$cms.Columns("B").ColumnWidth = 20.998697916666664
$cms.Columns("C").ColumnWidth = 18.166666666666664
$cms.Columns("E").ColumnWidth = 16.498697916666664

$cms.Range("E25").Select()

# --- ADCore sheet (sheet2) ---
$adcore = $wb.Worksheets.Item("ADCore")
$adcore.Activate()
$adcore.Range("C27").Value = "Banner 3/Blank"
$adcore.Rows("28").Insert()
$adcore.Range("A28").Value = "WXX_XX_XXX_XXXBN0400"
$adcore.Range("B28").Value = "W19_00_000_UPLD00 "
$adcore.Range("C28").Value = "Banner 4"
$adcore.Range("D28").Value = "WXX_XX_XXX_XXXBN0400"
$adcore.Range("G28").Value = "2018-full-mod"

$adcore.Range("C29").Select()

# Make sure CMS stays the active tab at the end (tabSelected=1 in original)
$cms.Activate()
